# "xen plan rates updated" -- refresh the Medgulf/Nextcare plan info row:
#   - residencies: drop the "AbuDhabi" half of the NE_Dubai/AbuDhabi pair
#   - conversion: new rate is 1 (was 3.6725)
#   - startDate: roll forward to 2025-01-01 (was 2024-10-31)
#   - currency: AED instead of USD
# D2 (startDate) also loses its special "text" font override, falling back
# to the sheet's default Arial font (same as the D1 header cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NE_Dubai/"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "2025-01-01"
$ws.Range("E2").Value = "AED"

# D2 reverts to the plain default font (matches D1's header formatting)
# instead of the bespoke font that used to be applied to it.
$ws.Range("D2").Font.Name = "Arial"
$ws.Range("D2").Font.FontStyle = "Regular"
$ws.Range("D2").Font.Size = 10

# Cursor/selection left on E3 after the edit.
$ws.Range("E3").Select() | Out-Null
